$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: single "transformer simple" result row, epochs label/value, and the
#     newly measured MIXED-dataset metrics (G3:J3) ---

# Intern "transformer simple" before "epochs" so the shared-string table order
# matches the authored workbook.
$ws.Range("B3").Value = "transformer simple"
$ws.Range("A2").Value = "epochs"
$ws.Range("A3").Value = 50

$ws.Range("C3").Value = 0.70389461498114902
$ws.Range("D3").Value = 0.661825089833181
$ws.Range("E3").Value = 13.240868748151801
$ws.Range("F3").Value = 0.26438746185434397

# G3 used to carry a bordered "empty" style; drop that left border so it
# matches its neighbours (H3:J3) before filling in the new MIXED-dataset values.
$ws.Range("G3").Borders.LineStyle = -4142

$ws.Range("G3").Value = 0.52420581224005502
$ws.Range("H3").Value = 0.82352536484161898
$ws.Range("I3").Value = 16.197356524641201
$ws.Range("J3").Value = 0.44115773940263903

# --- Remove now-obsolete model rows (CNN, LSTM, transformer) ---
$ws.Range("B4").ClearContents() | Out-Null
$ws.Range("B5").ClearContents() | Out-Null
$ws.Range("B6").ClearContents() | Out-Null

# --- Cosmetic updates matching the refreshed layout ---
$ws.Columns("B").ColumnWidth = 15.1666666666667
$ws.Range("C3:F3").Select() | Out-Null
